$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "data refreshed at" timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 07:05"

# --- India (row 14): refreshed case counters ---
$ws.Range("B14").Value = 101261
$ws.Range("C14").Value = 933
$ws.Range("E14").Value = 58864
$ws.Range("G14").Value = 8
$ws.Range("H14").Value = 3164

# --- Tailandia (row 73): refreshed case counters ---
$ws.Range("B73").Value = 3033
$ws.Range("C73").Value = 2
$ws.Range("E73").Value = 120

# --- Reorder Nueva Caledonia / Belice / Santa Lucia block (rows 195-197) ---
# New order: Belice, Santa Lucia, Nueva Caledonia (each row keeps its own data)
$ws.Range("A195").Value = "Belice"
$ws.Range("D195").Value = 16
$ws.Range("H195").Value = 2

$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Nueva Caledonia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# --- San Cristobal y Nieves (row 203): refreshed counters ---
$ws.Range("D203").Value = 15
$ws.Range("E203").Value = 0

# --- Reorder Seychelles / Groenlandia / Montserrat block (rows 209-211) ---
# New order: Montserrat, Seychelles, Groenlandia (each row keeps its own data)
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Groenlandia"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
